# Actualización 11 de Mayo - Mañana
# Adds 6 new "Rescatable" students to the "Rescatables" sheet. The two
# students that were already on the list (CABRERA / TORRES) stay on the
# list too, now at the bottom (rows 8 and 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Full data set for the sheet (NC, Paterno, Materno, Nombres, Nombre_Largo,
# Grupo, Reprobadas) once the new students are added.
$data = @(
    @(20330051920013, "FLORES",    "JUAREZ",    "LUIS ANGEL",       "GEOMETRÍA Y TRIGONOMETRÍA", "2AEM",  2),
    @(20330051920337, "BARRAGAN",  "VILLALBA",  "ADAN",             "GEOMETRÍA Y TRIGONOMETRÍA", "2APM",  2),
    @(20330051920381, "HERNANDEZ", "SANCHEZ",   "EDGAR DANIEL",     "GEOMETRÍA Y TRIGONOMETRÍA", "2APM",  2),
    @(20330051920373, "RICO",      "BAUTISTA",  "EDGAR RAMSES",     "GEOMETRÍA Y TRIGONOMETRÍA", "2ARHM", 2),
    @(20330051920313, "TORRES",    "VAZQUEZ",   "JOSELIN GUADALUPE","GEOMETRÍA Y TRIGONOMETRÍA", "2BLCM", 2),
    @(20330051920263, "CARRERA",   "ZAVALETA",  "ALFREDO",          "GEOMETRÍA Y TRIGONOMETRÍA", "2APM",  1),
    @(20330051920287, "CABRERA",   "RODRIGUEZ", "DANIEL",           "GEOMETRÍA Y TRIGONOMETRÍA", "2BLCM", 1),
    @(20330051920379, "TORRES",    "CARRASCO",  "ZULEICA RENATA",   "GEOMETRÍA Y TRIGONOMETRÍA", "2BLCM", 1)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
    $ws.Cells.Item($r, 6).Value2 = $row[5]
    $ws.Cells.Item($r, 7).Value2 = $row[6]
    $r = $r + 1
}
